$d = $word.ActiveDocument
$r = $d.Content

$found = $r.Find.Execute("uređenje", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r.Collapse(0)

    # First new run: ", forma "
    $run1Start = $r.Start
    $r.InsertAfter(", forma ")
    $run1End = $r.End
    $run1 = $d.Range($run1Start, $run1End)
    $run1.Font.NameAscii = "Calibri"
    $run1.Font.NameOther = "Calibri"
    $run1.Font.NameBi = "Calibri"
    $run1.Font.Color = 0

    # Second new run: "razlomka"
    $r.Collapse(0)
    $run2Start = $r.Start
    $r.InsertAfter("razlomka")
    $run2End = $r.End
    $run2 = $d.Range($run2Start, $run2End)
    $run2.Font.NameAscii = "Calibri"
    $run2.Font.NameOther = "Calibri"
    $run2.Font.NameBi = "Calibri"
    $run2.Font.Color = 0
}
